$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.892.99"
Set-TextValue $ws.Range("E2") "  +1.40%  "
Set-TextValue $ws.Range("D3") "3.414.84"
Set-TextValue $ws.Range("E3") "  +0.90%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "578.20"
Set-TextValue $ws.Range("E5") "  +1.24%  "
Set-TextValue $ws.Range("D6") "144.22"
Set-TextValue $ws.Range("E6") "  +2.25%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("E8") "  +0.11%  "
Set-TextValue $ws.Range("E9") "  -0.89%  "
Set-TextValue $ws.Range("E10") "  +0.48%  "
Set-TextValue $ws.Range("D11") "0.385"
Set-TextValue $ws.Range("E11") "  -0.57%  "
Set-TextValue $ws.Range("D12") "4.001.42"
Set-TextValue $ws.Range("E12") "  +0.98%  "
Set-TextValue $ws.Range("B13") "TRON"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.125"
Set-TextValue $ws.Range("E13") "  -0.84%  "
Set-TextValue $ws.Range("B14") "Avalanche"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D14") "28.36"
Set-TextValue $ws.Range("E14") "  +1.82%  "
Set-TextValue $ws.Range("D15") "3.411.40"
Set-TextValue $ws.Range("E15") "  +1.07%  "
Set-TextValue $ws.Range("E16") "  +0.06%  "
Set-TextValue $ws.Range("D17") "61.933.12"
Set-TextValue $ws.Range("E17") "  +1.32%  "
Set-TextValue $ws.Range("E18") "  +0.73%  "
Set-TextValue $ws.Range("D19") "13.98"
Set-TextValue $ws.Range("E19") "  +2.67%  "
Set-TextValue $ws.Range("E20") "  +2.87%  "
Set-TextValue $ws.Range("D21") "390.74"
Set-TextValue $ws.Range("E21") "  +2.24%  "
Set-TextValue $ws.Range("D22") "74.91"
Set-TextValue $ws.Range("E22") "  -1.32%  "
Set-TextValue $ws.Range("E23") "  +0.61%  "
Set-TextValue $ws.Range("E24") "  +0.28%  "
Set-TextValue $ws.Range("D25") "3.556.02"
Set-TextValue $ws.Range("E25") "  +1.03%  "
Set-TextValue $ws.Range("E26") "  -0.65%  "
Set-TextValue $ws.Range("E27") "  -2.65%  "
Set-TextValue $ws.Range("D28") "7.45"
Set-TextValue $ws.Range("E28") "  +2.42%  "
Set-TextValue $ws.Range("D29") "0.999"
Set-TextValue $ws.Range("E29") "  +0.02%  "
Set-TextValue $ws.Range("D30") "8.00"
Set-TextValue $ws.Range("E30") "  +0.57%  "
Set-TextValue $ws.Range("E31") "  +0.27%  "
Set-TextValue $ws.Range("E32") "  -0.01%  "
Set-TextValue $ws.Range("E33") "  +2.58%  "
Set-TextValue $ws.Range("D34") "23.57"
Set-TextValue $ws.Range("E34") "  +1.25%  "
Set-TextValue $ws.Range("D35") "5.25"
Set-TextValue $ws.Range("E35") "  +5.39%  "
Set-TextValue $ws.Range("D36") "6.97"
Set-TextValue $ws.Range("E36") "  +0.17%  "
Set-TextValue $ws.Range("D37") "167.64"
Set-TextValue $ws.Range("E38") "  +4.85%  "
Set-TextValue $ws.Range("D39") "3.446.90"
Set-TextValue $ws.Range("E39") "  +0.80%  "
Set-TextValue $ws.Range("D40") "28.43"
Set-TextValue $ws.Range("E40") "  +8.36%  "
Set-TextValue $ws.Range("D41") "0.0754"
Set-TextValue $ws.Range("E41") "  -1.77%  "
Set-TextValue $ws.Range("D42") "0.785"
Set-TextValue $ws.Range("E42") "  +0.62%  "
Set-TextValue $ws.Range("D43") "4.42"
Set-TextValue $ws.Range("E43") "  +1.33%  "
Set-TextValue $ws.Range("E44") "  +1.49%  "
Set-TextValue $ws.Range("E45") "  +4.16%  "
Set-TextValue $ws.Range("D46") "2.502.75"
Set-TextValue $ws.Range("E46") "  +2.24%  "
Set-TextValue $ws.Range("D47") "22.78"
Set-TextValue $ws.Range("E47") "  -0.89%  "
Set-TextValue $ws.Range("E48") "  -0.11%  "
Set-TextValue $ws.Range("E49") "  +0.04%  "
Set-TextValue $ws.Range("D50") "0.0262"
Set-TextValue $ws.Range("E50") "  +0.21%  "
Set-TextValue $ws.Range("D51") "2.08"
Set-TextValue $ws.Range("E51") "  -2.60%  "
